$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A2 text: remove leading space from " 1. CPB 279 PROJECTS"
$ws.Range("A2").Value = "1. CPB 279 PROJECTS"

# Row5: fill in B5, C5, D5 with 0 (previously blank)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

# Row13: E13 5 -> 7
$ws.Range("E13").Value = 7

# Column A width change (bestFit width increased due to longer text in A2)
$ws.Columns("A").ColumnWidth = 27.428571428571427

# sheetView zoom and selection
$ws.Application.ActiveWindow.Zoom = 123
$ws.Range("A16").Select()
